# "Connect unit with chassis"
# Adds datacenter/room/row placement columns (ЦОД / Помещение / Ряд) to the
# "Rack 1" sheet, fills them in for every equipment row, and fixes up a
# couple of pre-existing cells (H2, J54) that should hold the text "1"
# instead of being numeric/blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rack 1")

# ---- New header cells (row 1) ------------------------------------------
$ws.Cells.Item(1, 13).Value = "ЦОД"
$ws.Cells.Item(1, 14).Value = "Помещение"
$ws.Cells.Item(1, 15).Value = "Ряд"

# ---- Fix two pre-existing data cells ------------------------------------
# H2 used to be the number 1; it should now be the text "1".
$ws.Cells.Item(2, 8).Value = "1"
$ws.Range("H2").NumberFormat = "@"
$ws.Cells.Item(2, 8).Value = "1"

# J54 used to be empty; it should now hold the text "1".
$ws.Cells.Item(54, 10).Value = "1"
$ws.Range("J54").NumberFormat = "@"
$ws.Cells.Item(54, 10).Value = "1"

# ---- Populate ЦОД / Помещение / Ряд for every equipment row -------------
for ($r = 2; $r -le 66; $r++) {
    $ws.Cells.Item($r, 13).Value = "ЦОД-1"
    $ws.Cells.Item($r, 14).Value = 404
    $ws.Cells.Item($r, 15).Value = "B"
}

# ---- Column widths / bestfit --------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 66.85546875
$ws.Columns.Item(11).ColumnWidth = 51.7109375
$ws.Columns.Item(14).ColumnWidth = 12.140625
$ws.Columns.Item(15).ColumnWidth = 6.140625

# ---- Sheet view: scrolled to column D, selection on N49 -----------------
$ws.Activate()
$ws.Range("D1").Select()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("N49").Select()
